$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data.
# Rows 14/15 (Avalanche/WrappedBTC) and rows 22/23 (PEPE/Uniswap) swapped ranking positions.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel keeps
# them as text (matching the source data which stores all Price/Volume cells as strings).
$ws.Range("D2").Value = "93.208.23"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "3.126.60"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'242.91"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").Value = "'616.42"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'1.10"
$ws.Range("E7").Value = "  -3.81%  "
$ws.Range("D8").Value = "'0.403"
$ws.Range("E8").Value = "  +8.31%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "3.117.80"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "'0.736"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "'0.203"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'34.71"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "92.617.52"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'5.52"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "3.698.40"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "3.107.32"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "'3.73"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'14.79"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("D21").Value = "'5.84"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").Value = "'0.0000206"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'9.46"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").Value = "'450.19"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'5.82"
$ws.Range("E25").Value = "  -4.54%  "
$ws.Range("D26").Value = "'87.16"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").Value = "'11.81"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "3.280.21"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'0.136"
$ws.Range("E30").Value = "  +4.25%  "
$ws.Range("D31").Value = "'0.232"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").Value = "'9.32"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +12.62%  "
$ws.Range("D35").Value = "'8.08"
$ws.Range("E35").Value = "  +4.42%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E37").Value = "  +11.60%  "
$ws.Range("D38").Value = "'26.29"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").Value = "'486.26"
$ws.Range("E40").Value = "  -4.75%  "
$ws.Range("D41").Value = "'1.32"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").Value = "'3.52"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").Value = "'0.441"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").Value = "'23.11"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'161.26"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").Value = "'1.93"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "'0.698"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'0.0339"
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("D51").Value = "'4.44"
$ws.Range("E51").Value = "  -0.41%  "
Write-Output "Update complete"
